$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exceptions")

# New exception rows describing the JSON serialization work for the user profile feature
$ws.Range("A3").Value = "0x0001"
$ws.Range("C3").Value = "UserModel"
$ws.Range("B3").Value = "When signing in we found repeat emails in DB"

$ws.Range("A4").Value = "0x0002"
$ws.Range("B4").Value = "Failure displaying user profile because user is incomplete"
$ws.Range("C4").Value = "UserProfile.aspx"

# Move the selection to the next empty row, matching the post-edit state in Excel
$ws.Range("A5").Select()
